$d = $word.ActiveDocument

function Replace-First($old, $new, [bool]$wholeWord=$false) {
    # Search the whole document content for the first occurrence of $old
    # and overwrite that Range's .Text directly (bypasses autocorrect
    # smart-quote substitution that Find/Replace's ReplaceWith would apply).
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $wholeWord, $false, $false, $false, `
                                $true, 1, $false)
    if ($found) {
        $rng.Text = $new
        return $true
    }
    return $false
}

function Replace-All($old, $new, [bool]$wholeWord=$false) {
    $count = 0
    $rng = $d.Content
    $keepGoing = $true
    while ($keepGoing) {
        $found = $rng.Find.Execute($old, $true, $wholeWord, $false, $false, $false, `
                                    $true, 1, $false)
        if ($found) {
            $rng.Text = $new
            $count = $count + 1
            $rng.Collapse(0)
        } else {
            $keepGoing = $false
        }
    }
    return $count
}

# Title
Replace-First "Appendix 5: SWIFT Exit Interview Schedule" "Bylaag 5: SWIFT Uitgangsonderhoudskedule"

# Intro greeting
Replace-First "Hi! I just have three quick questions for you about your clinic visit today." "Hallo! Ek het drie vinnige vrae vir jou oor jou kliniekbesoek vandag."

# First (free-flow, trailing-space) occurrence of the three questions
Replace-First "Are you a parent or a caregiver of a child? " "Is jy 'n ouer of 'n versorger van 'n kind? "
Replace-First "Did you see a poster in the clinic today advertising a parenting programme?" "Het jy vandag 'n plakkaat in die kliniek gesien wat 'n ouerskapprogram adverteer?"
Replace-First "Did any of the nurses you saw today tell you about the programme?" "Het enige van die verpleegsters, wat jy vandag gesien het, jou oor die program vertel?"

# Interviewer instructions
Replace-First "Document for interviewers to keep track of and tally responses:" "Dokument vir onderhoudvoerders om antwoorde dop te hou en te tel:"
Replace-First "Site________________________      Clinic Name ____________________________   " "Instelling________________________      Kliniek Naam:____________________________   "
Replace-First "RA ______________________________       Date_______________________ " "RA ______________________________       Datum_______________________ "

# Second (table header, no trailing space) occurrence of the three questions
Replace-First "Are you a parent or a caregiver of a child?" "Is jy 'n ouer of 'n versorger van 'n kind?"
Replace-First "Did you see a poster in the clinic today advertising a parenting programme?" "Het jy vandag 'n plakkaat in die kliniek gesien wat 'n ouerskapprogram adverteer?"
Replace-First "Did any of the nurses you saw today tell you about the programme?" "Het enige van die verpleegsters, wat jy vandag gesien het, jou oor die program vertel?"

# Table answer cells: three Yes/No pairs
Replace-All "Yes" "Ja" $true
Replace-All "No" "Nee" $true
